$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Date (D) and Volume (M) values between row 2 and row 5
$ws.Range("D2").Value = 44981
$ws.Range("M2").Value = 30

$ws.Range("D5").Value = 44980
$ws.Range("M5").Value = 50
